# fix(publipostage): Correct status name
#
# Renames a handful of "statut_label" / "statut_name" values used by the
# mail-merge (publipostage) source sheet:
#   - "bleu"  -> "noir"
#   - "pas de résultat ni de publication"                    -> "pas de résultat postés ni publiés"
#   - "résultat et / ou publication posté"                   -> "résultat postés ou publiés"
#   - "résultat et / ou publication posté dans les 36 mois"  -> "résultat postés ou publiés dans les 36 mois"
#   - "résultat et / ou publication posté dans les 12 mois"  -> "résultat postés ou publiés dans les 12 mois"
#
# These strings appear many times throughout the used range (columns B and
# C), so a whole-cell Find & Replace across the sheet is used instead of
# touching individual cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlWhole = 1
$xlByRows = 1

# "bleu" -> "noir" (statut_label, column B)
$ws.Cells.Replace("bleu", "noir", $xlWhole, $xlByRows, $false, $false, $true, $true)

# statut_name values (column C) - replace the longer/more specific strings
# first so the shorter base phrase can't shadow them; whole-cell matching
# already makes the order irrelevant, but this keeps the intent obvious.
$ws.Cells.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois", $xlWhole, $xlByRows, $false, $false, $true, $true)
$ws.Cells.Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois", $xlWhole, $xlByRows, $false, $false, $true, $true)
$ws.Cells.Replace("résultat et / ou publication posté", "résultat postés ou publiés", $xlWhole, $xlByRows, $false, $false, $true, $true)
$ws.Cells.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés", $xlWhole, $xlByRows, $false, $false, $true, $true)
